# "update Key number and reports"
#
# Rename 7 header cells (row 1) to their "_py" suffixed variants and update
# the associated key-number counts in row 2 to the freshly recomputed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers (C1:I1)
$ws.Range("C1").Value = "orphaCoding_no_py"
$ws.Range("D1").Value = "orphaCase_no_py"
$ws.Range("E1").Value = "unambigous_rdCase_no_py"
$ws.Range("F1").Value = "rdCase_no_py"
$ws.Range("G1").Value = "case_no_py"
$ws.Range("H1").Value = "patient_no_py"
$ws.Range("I1").Value = "case_no_py_ipat"

# Row 2 key-number values (C2:H2); A2, B2, I2 are unchanged
$ws.Range("C2").Value = 92
$ws.Range("D2").Value = 90
$ws.Range("E2").Value = 80
$ws.Range("F2").Value = 99
$ws.Range("G2").Value = 997
$ws.Range("H2").Value = 950
